$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17-21 hold five occurrence records for the same taxon at the same
# locality; the edit re-shuffles which identifier/coordinate/phenology
# values land on which row (the underlying per-row metadata in the other
# columns is identical across these five rows).
#
# Capture the "identity" data (Id, Ost, Nord, and the optional
# Enhet/Ålder-Stadium/Kön/Metod/Bestämningsmetod block) for each of the
# five rows before any writes happen, then redistribute it according to
# the new row order.

$rows = 17..21

$idVal = @{}
$qVal  = @{}
$rVal  = @{}
$hasExtra = @{}
$kVal  = @{}

foreach ($r in $rows) {
    $idVal[$r] = $ws.Cells.Item($r, 1).Value2        # column A - Id
    $qVal[$r]  = $ws.Cells.Item($r, 17).Value2        # column Q - Ost
    $rVal[$r]  = $ws.Cells.Item($r, 18).Value2        # column R - Nord
    $kCell = $ws.Cells.Item($r, 11).Value2            # column K - Alder-Stadium
    if ($kCell -eq "blomning") {
        $hasExtra[$r] = $true
        $kVal[$r] = $kCell
    } else {
        $hasExtra[$r] = $false
    }
}

# new row -> source row the identity data now comes from
$sourceRow = @{ 17 = 21; 18 = 19; 19 = 17; 20 = 18; 21 = 20 }

foreach ($newR in $rows) {
    $src = $sourceRow[$newR]

    $ws.Cells.Item($newR, 1).Value = $idVal[$src]
    $ws.Cells.Item($newR, 17).Value = $qVal[$src]
    $ws.Cells.Item($newR, 18).Value = $rVal[$src]

    $jCell = $ws.Cells.Item($newR, 10)   # J - Enhet
    $kCellTgt = $ws.Cells.Item($newR, 11) # K - Alder-Stadium
    $lCell = $ws.Cells.Item($newR, 12)   # L - Kon
    $nCell = $ws.Cells.Item($newR, 14)   # N - Metod
    $afCell = $ws.Cells.Item($newR, 32)  # AF - Bestamningsmetod

    if ($hasExtra[$src]) {
        # Write the empty-but-typed marker cells (J, L, N, AF) and the
        # populated K cell ("blomning").
        $jCell.Value = "'"
        $jCell.Style = "Normal"

        $kCellTgt.Value = $kVal[$src]

        $lCell.Value = "'"
        $lCell.Style = "Normal"

        $nCell.Value = "'"
        $nCell.Style = "Normal"

        $afCell.Value = "'"
        $afCell.Style = "Normal"
    } else {
        # Make sure these cells are blank/absent on rows that shouldn't
        # carry the extra phenology block.
        $jCell.Value = ""
        $kCellTgt.Value = ""
        $lCell.Value = ""
        $nCell.Value = ""
        $afCell.Value = ""
    }
}
